$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 1
    4  = 2
    5  = 0
    6  = 1
    7  = 0
    8  = 2
    9  = 1
    10 = 3
    11 = 1
    12 = 1
    13 = 0
    14 = 2
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 2
    23 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
